# Excel COM-interop script: apply Linea 141 schedule update (scrape 20:13:01)
# Sheet1 = LP1912, Sheet2 = LP1912-215, Sheet3 = 6203-6173

$wb = $excel.ActiveWorkbook

$sheet1Data = @(
  @{Row=48; A="05:53:46"; B="07:31"; C="11_ETCHEVERRY"; D=98; E="LP1912"},
  @{Row=49; A="07:24:45"; B="07:31"; C="16_SANTA ANA"; D=7; E="LP1912"},
  @{Row=120; A="10:13:53"; B="10:34"; C="16_SANTA ANA"; D=21; E="LP1912"},
  @{Row=121; A="10:13:53"; B="10:34"; C="23_HERNANDEZ"; D=21; E="LP1912"},
  @{Row=162; A="10:13:53"; B="12:06"; C="14_ABASTO"; D=113; E="LP1912"},
  @{Row=163; A="10:52:37"; B="12:06"; C="10_OLMOS"; D=74; E="LP1912"},
  @{Row=164; A="10:13:53"; B="12:06"; C="16_P MOR-SANTA ANA"; D=113; E="LP1912"},
  @{Row=178; A="11:46:46"; B="12:34"; C="26_HERNANDEZ"; D=48; E="LP1912"},
  @{Row=179; A="11:17:39"; B="12:34"; C="11_ETCHEVERRY"; D=77; E="LP1912"},
  @{Row=193; A="10:52:37"; B="12:50"; C="15_ABASTO"; D=118; E="LP1912"},
  @{Row=194; A="12:50:41"; B="12:50"; C="16_SANTA ANA"; D=0; E="LP1912"},
  @{Row=220; A="12:01:11"; B="13:51"; C="215A_EL PATO"; D=110; E="LP1912"},
  @{Row=221; A="13:51:32"; B="13:51"; C="11_ETCHEVERRY"; D=0; E="LP1912"},
  @{Row=263; A="15:36:13"; B="15:37"; C="11_ETCHEVERRY"; D=1; E="LP1912"},
  @{Row=264; A="13:51:32"; B="15:37"; C="10_OLMOS"; D=106; E="LP1912"},
  @{Row=321; A="17:12:54"; B="17:20"; C="26_HERNANDEZ"; D=8; E="LP1912"},
  @{Row=322; A="16:45:34"; B="17:20"; C="16_SANTA ANA"; D=35; E="LP1912"},
  @{Row=335; A="15:59:02"; B="17:38"; C="17_ROMERO"; D=99; E="LP1912"},
  @{Row=336; A="16:45:34"; B="17:38"; C="27_EL RETIRO"; D=53; E="LP1912"},
  @{Row=364; A="18:13:12"; B="18:25"; C="14_ABASTO"; D=12; E="LP1912"},
  @{Row=365; A="17:50:30"; B="18:25"; C="26_HERNANDEZ"; D=35; E="LP1912"},
  @{Row=379; A="16:53:02"; B="18:48"; C="14X44_ABASTO"; D=115; E="LP1912"},
  @{Row=380; A="18:48:53"; B="18:48"; C="16_SANTA ANA"; D=0; E="LP1912"},
  @{Row=407; A="17:36:40"; B="19:30"; C="225_GOMEZ"; D=114; E="LP1912"},
  @{Row=408; A="18:35:28"; B="19:30"; C="16_SANTA ANA"; D=55; E="LP1912"},
  @{Row=431; A="20:13:01"; B="20:13"; C="16_SANTA ANA"; D=0; E="LP1912"},
  @{Row=432; A="20:13:01"; B="20:13"; C="17_ROMERO"; D=0; E="LP1912"},
  @{Row=433; A="18:48:53"; B="20:20"; C="26_HERNANDEZ"; D=92; E="LP1912"},
  @{Row=434; A="18:35:28"; B="20:21"; C="26_HERNANDEZ"; D=106; E="LP1912"},
  @{Row=435; A="18:48:53"; B="20:21"; C="11_ETCHEVERRY"; D=93; E="LP1912"},
  @{Row=436; A="18:35:28"; B="20:22"; C="11_ETCHEVERRY"; D=107; E="LP1912"},
  @{Row=437; A="19:56:05"; B="20:23"; C="16_SANTA ANA"; D=27; E="LP1912"},
  @{Row=438; A="18:35:28"; B="20:23"; C="215A_EL PATO"; D=108; E="LP1912"},
  @{Row=439; A="19:42:01"; B="20:24"; C="215A_EL PATO"; D=42; E="LP1912"},
  @{Row=440; A="20:13:01"; B="20:30"; C="14_ABASTO"; D=17; E="LP1912"},
  @{Row=441; A="18:48:53"; B="20:30"; C="225_GOMEZ"; D=102; E="LP1912"},
  @{Row=442; A="18:35:28"; B="20:31"; C="225_GOMEZ"; D=116; E="LP1912"},
  @{Row=443; A="19:42:01"; B="20:32"; C="14_ABASTO"; D=50; E="LP1912"},
  @{Row=444; A="19:56:05"; B="20:34"; C="14_ABASTO"; D=38; E="LP1912"},
  @{Row=445; A="20:13:01"; B="20:35"; C="16_SANTA ANA"; D=22; E="LP1912"},
  @{Row=446; A="19:42:01"; B="20:44"; C="11_ETCHEVERRY"; D=62; E="LP1912"},
  @{Row=447; A="19:16:50"; B="20:49"; C="11_ETCHEVERRY"; D=93; E="LP1912"},
  @{Row=448; A="19:16:50"; B="20:50"; C="14_ABASTO"; D=94; E="LP1912"},
  @{Row=449; A="19:42:01"; B="20:52"; C="15_ABASTO"; D=70; E="LP1912"},
  @{Row=450; A="19:42:01"; B="20:53"; C="23_HERNANDEZ"; D=71; E="LP1912"},
  @{Row=451; A="19:16:50"; B="20:55"; C="10_OLMOS"; D=99; E="LP1912"},
  @{Row=452; A="19:16:50"; B="20:56"; C="27_EL RETIRO"; D=100; E="LP1912"},
  @{Row=453; A="20:13:01"; B="20:56"; C="10_OLMOS"; D=43; E="LP1912"},
  @{Row=454; A="19:42:01"; B="20:57"; C="27_EL RETIRO"; D=75; E="LP1912"},
  @{Row=455; A="19:16:50"; B="21:04"; C="84_COLONIA URQUIZA-ESC 49"; D=108; E="LP1912"},
  @{Row=456; A="19:16:50"; B="21:07"; C="215B_EL PATO"; D=111; E="LP1912"},
  @{Row=457; A="19:42:01"; B="21:08"; C="215B_EL PATO"; D=86; E="LP1912"},
  @{Row=458; A="19:42:01"; B="21:21"; C="26_HERNANDEZ"; D=99; E="LP1912"},
  @{Row=459; A="20:13:01"; B="21:23"; C="15_ABASTO"; D=70; E="LP1912"},
  @{Row=460; A="19:42:01"; B="21:23"; C="10_OLMOS"; D=101; E="LP1912"},
  @{Row=461; A="19:42:01"; B="21:38"; C="14_ABASTO"; D=116; E="LP1912"},
  @{Row=462; A="19:42:01"; B="21:38"; C="17_ROMERO"; D=116; E="LP1912"},
  @{Row=463; A="19:56:05"; B="21:47"; C="215A_EL PATO"; D=111; E="LP1912"},
  @{Row=464; A="20:13:01"; B="22:08"; C="11_ETCHEVERRY"; D=115; E="LP1912"}
)

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 20:13:01"
$ws1.Range("A3").Value = "Total filas: 459"
foreach ($item in $sheet1Data) {
    $r = $item.Row
    $ws1.Range("A$r").Value = $item.A
    $ws1.Range("B$r").Value = $item.B
    $ws1.Range("C$r").Value = $item.C
    $ws1.Range("D$r").Value = $item.D
    $ws1.Range("E$r").Value = $item.E
}

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 20:13:01"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 20:13:01"
$ws3.Range("A3").Value = "Total filas: 57"
$ws3.Range("A62").Value = "20:13:01"
$ws3.Range("B62").Value = "22:05"
$ws3.Range("C62").Value = "215A_LA PLATA"
$ws3.Range("D62").Value = 112
$ws3.Range("E62").Value = "L6173"

Write-Output "Update applied successfully."
